$d = $word.ActiveDocument

# --- 1. Italicize the three package-name runs: "rstan", "rstanarm", "shiynstan" ---
# Use MatchWholeWord so "rstan" doesn't also match inside "rstanarm".

$rng = $d.Content
$null = $rng.Find.Execute("rstan", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Italic = 1

$rng = $d.Content
$null = $rng.Find.Execute("rstanarm", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Italic = 1

# --- 2. Move the "_GoBack" bookmark from the "Award Period" paragraph to
#        wrap the "shiynstan" run, and italicize that run too. ---

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$null = $rng.Find.Execute("shiynstan", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Italic = 1
$d.Bookmarks.Add("_GoBack", $rng)
